# Apply the DENSE FOG workbook update:
#  - Row 91: F91 (Like) 0 -> 1, M91 (Comments) 1 -> 2
#  - Add new Row 92 with a new Facebook post entry
#  - (new shared strings "DENSE FOG 2K21" and "10107942303217599" come
#    in automatically from the values written into D92 / E92)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- update existing row 91 ---
$ws.Range("F91").Value = 1
$ws.Range("M91").Value = 2

# --- append new row 92 ---
# Use raw serial numbers (rather than [DateTime]) so Excel doesn't invent a
# brand-new "general date" number format/style - we apply the same formats
# already used by the rest of column A/B explicitly below instead.
$ws.Range("A92").Value = 44199
$ws.Range("B92").Value = 0.49583333333333335
$ws.Range("C92").Value = "Friends"
$ws.Range("D92").Value = "DENSE FOG 2K21"
$ws.Range("E92").Value = "10107942303217599"
$ws.Range("F92").Value = 1
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1
$ws.Range("M92").Value = 3
$ws.Range("N92").Value = "10103925030244839"

# --- match formatting used for the other rows in these columns ---
$ws.Range("A92").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("B92").NumberFormat = "h:mm:ss;@"
$ws.Range("E92").NumberFormat = "@"
$ws.Range("N92").NumberFormat = "@"

# --- keep the view roughly where Excel would leave it after this edit ---
$ws.Application.ActiveWindow.ScrollRow = 86
$ws.Range("E91").Select()
